$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Replace "int" with "uuid" in the data-type column for specific fields
# (these cells currently hold the shared string "int")
$cells = @("E3", "I3", "O3", "S3", "O4", "I11", "S11", "I12")
foreach ($addr in $cells) {
    $ws.Range($addr).Value = "uuid"
}

# Update the active selection to S12
$ws.Range("S12").Select()
